$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update LOOKUP threshold table (D34:D49): shift bucket boundaries down by
# ~50 so that export_peak is computed correctly as min(import_peak,
# export_total) instead of being offset.
$values = @(1,201,251,301,351,401,451,501,551,601,651,701,751,801,851,901)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 34 + $i
    $ws.Cells.Item($row, 4).Value2 = $values[$i]
}

# Fix the LOOKUP formulas that previously added an erroneous +50 offset
$ws.Range("E17").Formula = "=LOOKUP(E2,D34:D49,E34:E49)*E2"
$ws.Range("E28").Formula = "=LOOKUP(E2,D34:D49,E34:E49)*-E5"

# Restore cursor/selection position as left by the author
$ws.Range("E29").Select()
